$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.437.65'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '3.848.16'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.02'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.02'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("D7").Value = '3.848.39'
$ws.Range("E7").Value = '  -1.19%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  -2.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.48'
$ws.Range("E11").Value = '  +1.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  -2.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000267'
$ws.Range("E13").Value = '  +4.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.09'
$ws.Range("E14").Value = '  -3.02%  '
$ws.Range("D15").Value = '4.494.54'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").Value = '3.851.02'
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").Value = '68.509.08'
$ws.Range("E17").Value = '  -1.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.52'
$ws.Range("E18").Value = '  -1.54%  '
$ws.Range("E19").Value = '  -3.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.20'
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '471.00'
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("E24").Value = '  -3.83%  '
$ws.Range("E25").Value = '  -2.15%  '
$ws.Range("E26").Value = '  -2.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.11'
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("E28").Value = '  +0.94%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").Value = '3.999.01'
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.48'
$ws.Range("E33").Value = '  -1.57%  '
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("E35").Value = '  -2.94%  '
$ws.Range("D36").Value = '3.815.03'
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.78'
$ws.Range("E37").Value = '  +11.14%  '
$ws.Range("E38").Value = '  -2.23%  '
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("E40").Value = '  -2.56%  '
$ws.Range("E41").Value = '  -2.67%  '
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("E43").Value = '  -3.26%  '
$ws.Range("E44").Value = '  -4.18%  '
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '418.33'
$ws.Range("E46").Value = '  -4.36%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000294'
$ws.Range("E48").Value = '  +5.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '47.00'
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.04'
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0358'
$ws.Range("E51").Value = '  -2.01%  '
